$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G4" = 2.7
    "H4" = 2.7
    "I4" = 3.2
    "R4" = 2
    "S4" = 1.75
    "T4" = 7
    "W4" = 26
    "X4" = 26
    "AB4" = 15
    "AG4" = 12
    "G8" = 1.91
    "I8" = 4.2
    "R8" = 2.1
    "S8" = 1.67
    "W8" = 15
    "AE8" = 9.5
    "AF8" = 21
    "J11" = 1.13
    "K11" = 6
    "G12" = 2.77
    "H12" = 3
    "L12" = 1.45
    "M12" = 2.4
    "N12" = 2.27
    "O12" = 1.5
    "P12" = 1.5
    "Q12" = 2.27
    "R12" = 1.98
    "S12" = 1.65
    "T12" = 7.1
    "V12" = 10.75
    "X12" = 28
    "Y12" = 45
    "Z12" = 6.9
    "AB12" = 17
    "AC12" = 100
    "AE12" = 6.5
    "AG12" = 10.25
    "AH12" = 28
    "AI12" = 26
    "AJ12" = 45
    "H13" = 5.1
    "I13" = 13.5
    "M13" = 3.4
    "R13" = 2.32
    "T13" = 6.1
    "V13" = 9.5
    "X13" = 11.75
    "Z13" = 10.75
    "AA13" = 11
    "AE13" = 28
    "AF13" = 120
    "AG13" = 50
    "AH13" = 700
    "AI13" = 300
    "K17" = 13
    "N17" = 1.83
    "O17" = 2.03
    "Z20" = 7.1
    "AB20" = 17.5
    "AE20" = 7.3
    "G23" = 2.35
    "H23" = 2.88
    "I23" = 3
    "J23" = 1.11
    "K23" = 6.5
    "T23" = 6.5
    "V23" = 10
    "W23" = 23
    "AE23" = 8
    "AF23" = 15
    "AH23" = 34
    "G24" = 1.12
    "I24" = 26
    "R24" = 2.47
    "V24" = 11
    "Z24" = 13
    "AB24" = 45
    "AE24" = 65
    "AF24" = 350
    "AI24" = 700
    "AJ24" = 350
    "H26" = 3.5
    "I26" = 4
    "R26" = 1.78
    "S26" = 1.83
    "T26" = 7
    "U26" = 8.5
    "W26" = 15
    "Z26" = 9.75
    "AA26" = 6.8
    "AB26" = 15.5
    "AC26" = 75
    "AE26" = 11
    "AF26" = 22
    "AH26" = 60
    "AI26" = 40
    "I27" = 4.1
    "J27" = 1.03
    "L27" = 1.22
    "N27" = 1.8
    "T27" = 8
    "U27" = 9
    "Y27" = 23
    "AA27" = 7
    "AC27" = 41
    "G60" = 2.15
    "H60" = 3.25
    "I60" = 3.15
    "L60" = 1.4
    "M60" = 2.5
    "N60" = 2.18
    "O60" = 1.53
    "P60" = 1.5
    "Q60" = 2.27
    "R60" = 1.93
    "S60" = 1.7
    "T60" = 6.2
    "U60" = 9.25
    "V60" = 9.25
    "W60" = 20
    "X60" = 20
    "Y60" = 37
    "Z60" = 7.7
    "AA60" = 6.3
    "AB60" = 17.5
    "AC60" = 100
    "AE60" = 7.9
    "AF60" = 15
    "AG60" = 11.75
    "AH60" = 40
    "AI60" = 32
    "AJ60" = 50
    "J28" = 1.05
    "L28" = 1.47
    "J29" = 1.03
    "K29" = 12
    "L29" = 1.19
    "M29" = 4
    "N29" = 1.8
    "O29" = 2
    "J30" = 1.04
    "L30" = 1.22
    "N30" = 1.87
    "O30" = 1.87
    "J34" = 1.05
    "K34" = 11
    "G38" = 2.87
    "I38" = 2.32
    "M38" = 2.95
    "N38" = 2.02
    "O38" = 1.62
    "P38" = 1.38
    "Q38" = 2.47
    "T38" = 6.9
    "U38" = 11.75
    "V38" = 9
    "W38" = 28
    "X38" = 21
    "Y38" = 29
    "Z38" = 8
    "AB38" = 11.5
    "AE38" = 6.4
    "AF38" = 9.5
    "AG38" = 7.7
    "AH38" = 19.5
    "AJ38" = 23
    "G39" = 1.33
    "R39" = 2.62
    "S39" = 1.41
    "J40" = 1.08
    "K40" = 8
    "G41" = 1.6
    "H41" = 3.75
    "I41" = 6
    "J41" = 1.07
    "K41" = 9
    "R41" = 2
    "S41" = 1.75
    "U41" = 7
    "AD41" = 401
    "AE41" = 13
    "G46" = 2.1
    "I46" = 3.25
    "M46" = 3.2
    "T46" = 8.75
    "V46" = 8.5
    "W46" = 21
    "X46" = 16
    "Y46" = 23
    "Z46" = 10.5
    "AB46" = 12.5
    "AE46" = 10.5
    "AF46" = 18
    "AG46" = 11.25
    "AH46" = 45
    "AI46" = 29
    "AJ46" = 32
    "G50" = 2.05
    "H50" = 3.6
    "I50" = 3.3
    "L50" = 1.11
    "M50" = 6.5
    "Y50" = 17
    "AA50" = 8
    "AD50" = 67
    "AE50" = 19
    "AG50" = 13
    "G52" = 3.7
    "J52" = 1.01
    "L52" = 1.11
    "J53" = 1.02
    "L53" = 1.17
    "J54" = 1.03
    "L54" = 1.22
    "J55" = 1.02
    "L55" = 1.15
    "T55" = 11
    "G56" = 4
    "H56" = 3.6
    "I56" = 1.83
    "J56" = 1.03
    "L56" = 1.19
    "N56" = 1.73
    "O56" = 2.08
    "Z56" = 13
    "AF56" = 9.5
    "AH56" = 15
    "J57" = 1.03
    "L57" = 1.19
    "H62" = 5
    "K62" = 21
    "L62" = 1.11
    "M62" = 6.5
    "N62" = 1.4
    "O62" = 2.88
    "R62" = 1.67
    "S62" = 2.1
    "X62" = 10
    "AD62" = 350
    "AE62" = 29
    "AH62" = 101
    "N63" = 1.57
    "N64" = 1.77
    "O67" = 1.63
    "G68" = 1.85
    "I68" = 3.75
    "L68" = 1.29
    "M68" = 3.5
    "N68" = 1.93
    "O68" = 1.93
    "W68" = 17
    "AD68" = 251
    "AE68" = 11
    "G71" = 2.4
    "I71" = 3
    "K71" = 8.5
    "L71" = 1.33
    "M71" = 3.25
    "N71" = 2.1
    "O71" = 1.7
    "U71" = 11
    "Z71" = 8.5
    "AD71" = 251
    "AE71" = 9
    "AF71" = 15
    "AG71" = 12
    "AH71" = 34
    "G73" = 1.39
    "H73" = 4.3
    "I73" = 7.1
    "R73" = 1.91
    "S73" = 1.7
    "T73" = 6.8
    "U73" = 6.5
    "W73" = 9
    "Y73" = 28
    "AB73" = 20
    "AC73" = 100
    "AE73" = 18
    "AG73" = 23
    "G74" = 2.62
    "H74" = 3.4
    "I74" = 2.42
    "V74" = 9.75
    "X74" = 20
    "Y74" = 25
    "AG74" = 9.25
    "AI74" = 18.5
    "AJ74" = 24
    "K76" = 7.3
    "L76" = 1.31
    "M76" = 3.15
    "N76" = 1.91
    "P76" = 1.39
    "Q76" = 2.77
    "V76" = 8.5
    "Z76" = 7.3
    "AE76" = 13.5
    "G77" = 1.8
    "H77" = 3.3
    "I77" = 4.1
    "J77" = 1.07
    "K77" = 6.8
    "Q77" = 2.72
    "R77" = 1.87
    "S77" = 1.83
    "T77" = 6.4
    "U77" = 8.25
    "V77" = 8.25
    "W77" = 15
    "X77" = 15
    "Y77" = 29
    "Z77" = 6.8
    "AB77" = 15.5
    "AC77" = 80
    "AD77" = 700
    "AE77" = 11
    "AF77" = 23
    "AG77" = 14
    "AI77" = 40
    "AJ77" = 50
    "K79" = 6.4
    "L79" = 1.38
    "M79" = 2.8
    "N79" = 2.12
    "O79" = 1.65
    "P79" = 1.5
    "Q79" = 2.42
    "R79" = 1.85
    "S79" = 1.85
    "T79" = 8.25
    "Z79" = 6.4
    "AB79" = 14.5
    "AC79" = 75
    "AD79" = 700
    "AE79" = 7
    "AF79" = 11
    "AG79" = 9.5
    "AI79" = 21
    "AJ79" = 35
    "L80" = 1.33
    "M80" = 3.05
    "N80" = 1.98
    "O80" = 1.75
    "R80" = 1.85
    "T80" = 6.8
    "AG80" = 12.5
    "J81" = 1.04
    "K81" = 13
    "L81" = 1.2
    "M81" = 4.33
    "G85" = 2.3
    "I85" = 3
    "P85" = 1.4
    "Q85" = 2.75
    "U85" = 11
    "AD85" = 251
    "AG85" = 11
    "G86" = 2.35
    "I86" = 3.5
    "J86" = 1.11
    "K86" = 6.5
    "N86" = 2.5
    "Q86" = 2.38
    "Y86" = 41
    "AB86" = 17
    "AC86" = 67
    "AE86" = 8.5
    "AF86" = 15
    "AI68" = 29
    "AJ68" = 34
    "AI73" = 90
    "AJ73" = 80
    "L77" = 1.34
    "M77" = 3
    "N77" = 2
    "O77" = 1.72
    "P77" = 1.4
    "O86" = 1.5
    "P86" = 1.53
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
